# Apply "Trade #25 closed" update to the live trading results workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet - update aggregated metrics
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.88   # Current Capital
$summary.Range("B4").Value = -0.12     # Total P&L $
$summary.Range("B5").Value = -0.1      # Total P&L %
$summary.Range("B6").Value = 25        # Total Trades
$summary.Range("B8").Value = 12        # Losing Trades
$summary.Range("B9").Value = 28        # Win Rate %

# ---------------------------------------------------------------------------
# 2. Strategy Status sheet - update MarketMaking strategy row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.88   # Capital
$status.Range("D4").Value = 25      # Trades
$status.Range("E4").Value = -0.12   # P&L $
$status.Range("F4").Value = -0.12   # P&L %
$status.Range("G4").Value = 28      # Win Rate %

# ---------------------------------------------------------------------------
# 3. Append the newly closed trade (#25) to both "All Trades" and
#    "MarketMaking" sheets - they carry the same trade log data.
# ---------------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A26").Value = 25

    $ws.Range("B26").NumberFormat = "@"
    $ws.Range("B26").Value = "2026-02-17"
    $ws.Range("B26").Style = "Normal"

    $ws.Range("C26").NumberFormat = "@"
    $ws.Range("C26").Value = "15:19:21"
    $ws.Range("C26").Style = "Normal"

    $ws.Range("D26").Value = "MarketMaking"
    $ws.Range("E26").Value = "UP"
    $ws.Range("F26").Value = 0.9
    $ws.Range("G26").Value = 0.79
    $ws.Range("H26").Value = "CLOSED"
    $ws.Range("I26").Value = -12.2222
    $ws.Range("J26").Value = -0.11
    $ws.Range("K26").Value = 99.88
    $ws.Range("L26").Value = 0
    $ws.Range("M26").Value = 0
    $ws.Range("N26").Value = 0.6
    $ws.Range("O26").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P26").Value = "early_exit"
    $ws.Range("Q26").Value = 0.18
}
